$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Pepino ensalada needs to be inserted as the
# new row 270, pushing the existing rows 270-353 down to 271-354.
$ws.Rows("270:270").Insert()

$ws.Cells.Item(270, 1).Value2 = 7
$ws.Cells.Item(270, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(270, 3).Value2 = "Ñuble"
$ws.Cells.Item(270, 4).Value2 = 45120
$ws.Cells.Item(270, 5).Value2 = 16
$ws.Cells.Item(270, 6).Value2 = 100112043
$ws.Cells.Item(270, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(270, 8).Value2 = "Sin especificar"
$ws.Cells.Item(270, 9).Value2 = "Primera"
$ws.Cells.Item(270, 10).Value2 = 100
$ws.Cells.Item(270, 11).Value2 = 14000
$ws.Cells.Item(270, 12).Value2 = 14000
$ws.Cells.Item(270, 13).Value2 = 14000
$ws.Cells.Item(270, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(270, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(270, 16).Value2 = 233
$ws.Cells.Item(270, 17).Value2 = 60
$ws.Cells.Item(270, 18).Value2 = "Hortaliza"
